# Regenerate the "K" (strikeouts) column (column G) values for rows 2-46
# of Sheet1, replacing the previous "Strike#" counts with freshly
# calculated K values, per commit message:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values, keyed by row number (2-46)
$kValues = @{
    2  = 3
    3  = 1
    4  = 5
    5  = 0
    6  = 2
    7  = 4
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 2
    13 = 0
    14 = 3
    15 = 2
    16 = 0
    17 = 1
    18 = 2
    19 = 0
    20 = 2
    21 = 0
    22 = 4
    23 = 0
    24 = 1
    25 = 0
    26 = 2
    27 = 1
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 3
    36 = 1
    37 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 1
    42 = 1
    43 = 1
    44 = 1
    45 = 1
    46 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
